$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix row 349: the "NA" value was incorrectly stored in the Abundance column (E).
# Move it to the Observations column (F) with clearer English wording, and blank
# out the Abundance cell. Doing this before the header rename keeps the shared
# string table ordering aligned with the canonical output.
$ws.Range("E349").Value = ""
$ws.Range("F349").Value = "Dato no apuntado"

# Translate / standardize the header row (also renames the Excel Table columns).
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Field"
$ws.Range("C1").Value = "Treatment"
$ws.Range("E1").Value = "Abundance"
$ws.Range("F1").Value = "Observations"
$ws.Range("D1").Value = "Repeat"

# Rename the worksheet itself (English name).
$ws.Name = "Lisso larvae"

# Resize columns now that the (shorter) English headers no longer need the
# previous "best fit" widths.
$ws.Columns.Item(1).ColumnWidth = 10.33
$ws.Columns.Item(2).ColumnWidth = 7.67
$ws.Columns.Item(3).ColumnWidth = 11.5
$ws.Columns.Item(4).ColumnWidth = 11.17
$ws.Columns.Item(5).ColumnWidth = 11.83

# Reset the view: scroll back to the top and select F10.
$null = $ws.Range("F10").Select()
